# Apply updated market-price / profit values per sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 4983.5
$ws.Range("J7").Value = 4983.5
$ws.Range("L7").Value = 4983.5
$ws.Range("N7").Value = -5207.5
$ws.Range("H10").Value = 39000
$ws.Range("J10").Value = 39000
$ws.Range("L10").Value = 39000
$ws.Range("N10").Value = -39586
$ws.Range("H13").Value = 2383.3333
$ws.Range("J13").Value = 3500
$ws.Range("L13").Value = 3500
$ws.Range("N13").Value = -3838
$ws.Range("H14").Value = 4983.5
$ws.Range("J14").Value = 4983.5
$ws.Range("L14").Value = 4983.5
$ws.Range("N14").Value = -5365.5
$ws.Range("H16").Value = 9900
$ws.Range("J16").Value = 9900
$ws.Range("L16").Value = 9900
$ws.Range("N16").Value = -10360
$ws.Range("H34").Value = 4999.5
$ws.Range("I34").Value = 4999.5
$ws.Range("K34").Value = 4999.5
$ws.Range("M34").Value = -4796.5
$ws.Range("H36").Value = 4999.5
$ws.Range("I36").Value = 4999.5
$ws.Range("K36").Value = 4999.5
$ws.Range("M36").Value = -4284.5
$ws.Range("H38").Value = 1972.5714
$ws.Range("I38").Value = 286.75
$ws.Range("J38").Value = 2646.9
$ws.Range("K38").Value = 860.25
$ws.Range("L38").Value = 7940.700000000001
$ws.Range("M38").Value = -488.25
$ws.Range("N38").Value = -8684.700000000001
$ws.Range("H82").Value = 1946.25
$ws.Range("I82").Value = 1946.25
$ws.Range("K82").Value = 5838.75
$ws.Range("M82").Value = -5432.75
$ws.Range("H85").Value = 1946.25
$ws.Range("I85").Value = 1946.25
$ws.Range("K85").Value = 5838.75
$ws.Range("M85").Value = -4434.75
$ws.Range("H116").Value = 3433.7778
$ws.Range("I116").Value = 3141
$ws.Range("K116").Value = 3141
$ws.Range("M116").Value = 301
$ws.Range("H132").Value = 1268.1459
$ws.Range("I132").Value = 1074.3636
$ws.Range("K132").Value = 3223.0908
$ws.Range("M132").Value = -693.0907999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2186.9443
$ws.Range("I2").Value = 1522.8125
$ws.Range("K2").Value = 1522.8125
$ws.Range("M2").Value = -1409.8125
$ws.Range("H32").Value = 17156
$ws.Range("I32").Value = 17156
$ws.Range("K32").Value = 17156
$ws.Range("M32").Value = -16869
$ws.Range("H37").Value = 4750
$ws.Range("H44").Value = 25500
$ws.Range("H55").Value = 22260
$ws.Range("J55").Value = 19680
$ws.Range("L55").Value = 19680
$ws.Range("N55").Value = -20310
$ws.Range("H116").Value = 2186.9443
$ws.Range("I116").Value = 1522.8125
$ws.Range("K116").Value = 1522.8125
$ws.Range("M116").Value = 771.1875
$ws.Range("H137").Value = 89999.836
$ws.Range("J137").Value = 89999.836
$ws.Range("L137").Value = 89999.836
$ws.Range("N137").Value = -100199.836

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2186.9443
$ws.Range("I3").Value = 1522.8125
$ws.Range("K3").Value = 1522.8125
$ws.Range("M3").Value = -1408.8125
$ws.Range("H82").Value = 24532
$ws.Range("J82").Value = 36886.668
$ws.Range("L82").Value = 36886.668
$ws.Range("N82").Value = -37652.668
$ws.Range("H85").Value = 24532
$ws.Range("J85").Value = 36886.668
$ws.Range("L85").Value = 36886.668
$ws.Range("N85").Value = -39538.668
$ws.Range("H97").Value = 15925.333
$ws.Range("I97").Value = 11343.25
$ws.Range("K97").Value = 11343.25
$ws.Range("M97").Value = -10352.25
$ws.Range("H105").Value = 2954.8462
$ws.Range("I105").Value = 2970.2727
$ws.Range("K105").Value = 2970.2727
$ws.Range("M105").Value = -1223.2727

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 55561590
$ws.Range("J31").Value = 11699.857
$ws.Range("L31").Value = 11699.857
$ws.Range("N31").Value = -12289.857
$ws.Range("H34").Value = 55561590
$ws.Range("J34").Value = 11699.857
$ws.Range("L34").Value = 11699.857
$ws.Range("N34").Value = -12103.857
$ws.Range("H58").Value = 6595.6
$ws.Range("I58").Value = 2635.6365
$ws.Range("K58").Value = 2635.6365
$ws.Range("M58").Value = -2432.6365
$ws.Range("H68").Value = 42864.5
$ws.Range("I68").Value = 42800
$ws.Range("K68").Value = 42800
$ws.Range("M68").Value = -42051
$ws.Range("H71").Value = 42864.5
$ws.Range("I71").Value = 42800
$ws.Range("K71").Value = 128400
$ws.Range("M71").Value = -124656
$ws.Range("H134").Value = 8071.875
$ws.Range("I134").Value = 7500
$ws.Range("K134").Value = 22500
$ws.Range("M134").Value = -19965
$ws.Range("H136").Value = 6595.6
$ws.Range("I136").Value = 2635.6365
$ws.Range("K136").Value = 7906.9095
$ws.Range("M136").Value = -5356.9095

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 999.6667
$ws.Range("I5").Value = 999
$ws.Range("K5").Value = 2997
$ws.Range("M5").Value = -2885
$ws.Range("H41").Value = 99
$ws.Range("J41").Value = 99
$ws.Range("L41").Value = 297
$ws.Range("N41").Value = -973
$ws.Range("H68").Value = 55559744
$ws.Range("I68").Value = 166670670
$ws.Range("J68").Value = 4282.8335
$ws.Range("K68").Value = 500012010
$ws.Range("L68").Value = 12848.5005
$ws.Range("M68").Value = -500011199
$ws.Range("N68").Value = -14470.5005
$ws.Range("H71").Value = 55559744
$ws.Range("I71").Value = 166670670
$ws.Range("J71").Value = 4282.8335
$ws.Range("K71").Value = 1500036030
$ws.Range("L71").Value = 38545.5015
$ws.Range("M71").Value = -1500031974
$ws.Range("N71").Value = -46657.5015
$ws.Range("H81").Value = 172285.42
$ws.Range("J81").Value = 172285.42
$ws.Range("L81").Value = 516856.26
$ws.Range("N81").Value = -519102.26
$ws.Range("H84").Value = 172285.42
$ws.Range("J84").Value = 172285.42
$ws.Range("L84").Value = 1550568.78
$ws.Range("N84").Value = -1561800.78
$ws.Range("H107").Value = 270618.94
$ws.Range("J107").Value = 500335.1
$ws.Range("L107").Value = 1501005.3
$ws.Range("N107").Value = -1504845.3
$ws.Range("H109").Value = 6081.75
$ws.Range("I109").Value = 2513.5
$ws.Range("K109").Value = 7540.5
$ws.Range("M109").Value = -6500.5
$ws.Range("H135").Value = 999.6667
$ws.Range("I135").Value = 999
$ws.Range("K135").Value = 8991
$ws.Range("M135").Value = -6456
$ws.Range("H137").Value = 8556.299999999999
$ws.Range("I137").Value = 6472.75
$ws.Range("J137").Value = 9945.333000000001
$ws.Range("K137").Value = 19418.25
$ws.Range("L137").Value = 29835.999
$ws.Range("M137").Value = -14318.25
$ws.Range("N137").Value = -40035.999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 35420
$ws.Range("J46").Value = 35420
$ws.Range("L46").Value = 35420
$ws.Range("N46").Value = -35732
$ws.Range("H121").Value = 48243.5
$ws.Range("J121").Value = 48243.5
$ws.Range("L121").Value = 48243.5
$ws.Range("N121").Value = -51737.5
$ws.Range("H122").Value = 2347.077
$ws.Range("I122").Value = 2391.3
$ws.Range("K122").Value = 7173.900000000001
$ws.Range("M122").Value = -4723.900000000001
$ws.Range("H123").Value = 54813
$ws.Range("J123").Value = 54813
$ws.Range("L123").Value = 54813
$ws.Range("N123").Value = -59713
$ws.Range("H132").Value = 1892.2059
$ws.Range("I132").Value = 1628.4482
$ws.Range("K132").Value = 4885.3446
$ws.Range("M132").Value = -2355.3446
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("H136").Value = 21625.812
$ws.Range("J136").Value = 21625.812
$ws.Range("L136").Value = 64877.436
$ws.Range("N136").Value = -69977.436

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H39").Value = 24995
$ws.Range("J39").Value = 24995
$ws.Range("L39").Value = 24995
$ws.Range("N39").Value = -25915
$ws.Range("H40").Value = 3144.5454
$ws.Range("I40").Value = 3198.125
$ws.Range("K40").Value = 3198.125
$ws.Range("M40").Value = -3062.125
$ws.Range("H136").Value = 100003700
$ws.Range("I136").Value = 58826700
$ws.Range("J136").Value = 333340000
$ws.Range("K136").Value = 176480100
$ws.Range("L136").Value = 1000020000
$ws.Range("M136").Value = -176477550
$ws.Range("N136").Value = -1000025100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H99").Value = 39914.4
$ws.Range("J99").Value = 39924.5
$ws.Range("L99").Value = 39924.5
$ws.Range("N99").Value = -45914.5
$ws.Range("H107").Value = 1123.7333
$ws.Range("I107").Value = 1059.9048
$ws.Range("J107").Value = 1272.6666
$ws.Range("K107").Value = 3179.7144
$ws.Range("L107").Value = 3817.9998
$ws.Range("M107").Value = -1259.7144
$ws.Range("N107").Value = -7657.9998
$ws.Range("H126").Value = 3622.8928
$ws.Range("I126").Value = 3349.6086
$ws.Range("K126").Value = 10048.8258
$ws.Range("M126").Value = -7578.825800000001
